$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the first "Textkoerper" (Body Text) paragraph of the report body --
# the one starting with "Im Rahmen dieses Projektberichts ...". We need to
# insert three new paragraphs right before it:
#   1. An "Abstract" styled paragraph (Zusammenfassung ...)
#   2. A "key words" styled paragraph (Schlagwoerter ...)
#   3. A "Heading 1" styled paragraph ("Einleitung")
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Im Rahmen dieses Projektberichts*") {
        $target = $i
        break
    }
}

if ($target -eq $null) {
    throw "could not locate target paragraph"
}

# Insert 3 empty paragraphs right before the target paragraph.
$targetRange = $d.Paragraphs.Item($target).Range
$targetRange.InsertParagraphBefore()
$targetRange.InsertParagraphBefore()
$targetRange.InsertParagraphBefore()

# After the three inserts, the target paragraph moved down by 3; the three
# new (still empty, "Body Text" styled) paragraphs are immediately before it.
$abstractPara = $d.Paragraphs.Item($target)
$keywordsPara = $d.Paragraphs.Item($target + 1)
$headingPara  = $d.Paragraphs.Item($target + 2)

# ---------------------------------------------------------------------------
# 1. Abstract paragraph
# ---------------------------------------------------------------------------
$abstractPara.Style = $d.Styles.Item("Abstract")
$ip = $d.Range($abstractPara.Range.Start, $abstractPara.Range.Start)

$ip.InsertAfter("Zusammenfassung")
$ip.Font.Italic = 1
$ip.Font.NameFarEast = "MS Mincho"
$ip.LanguageID = "de-DE"
$ip.Collapse(0)

$ip.InsertAfter("—")
$ip.Font.Italic = 0
$ip.Font.NameFarEast = "MS Mincho"
$ip.LanguageID = "de-DE"
$ip.Collapse(0)

$ip.InsertAfter(" Dieser Bericht beschreibt das Erstellen einer Deep Learning Architektur, welche die Handzeichen „Stein, Schere und Papier“ erkennen soll. Verwendet werden dabei die Datensätze von „Julien de la ")
$ip.Font.NameFarEast = "MS Mincho"
$ip.LanguageID = "de-DE"
$ip.Collapse(0)

$ip.InsertAfter("Bruère-Terreault")
$ip.Font.NameFarEast = "MS Mincho"
$ip.LanguageID = "de-DE"
$ip.Collapse(0)

$ip.InsertAfter("“ [1], sowie ein eigens erstellter Datensatz. ")
$ip.Font.NameFarEast = "MS Mincho"
$ip.LanguageID = "de-DE"
$ip.Collapse(0)

# ---------------------------------------------------------------------------
# 2. Keywords paragraph
# ---------------------------------------------------------------------------
$keywordsPara.Style = $d.Styles.Item("key words")
$ip2 = $d.Range($keywordsPara.Range.Start, $keywordsPara.Range.Start)

$ip2.InsertAfter("Schlagwörter—Deep Learning, Bilderkennung, Künstliche Intelligenz, Klassifizierung")
$ip2.Font.NameFarEast = "MS Mincho"
$ip2.LanguageID = "de-DE"
$ip2.Collapse(0)

# ---------------------------------------------------------------------------
# 3. "Einleitung" heading paragraph
# ---------------------------------------------------------------------------
$headingPara.Style = $d.Styles.Item("Heading 1")
$ip3 = $d.Range($headingPara.Range.Start, $headingPara.Range.Start)
$ip3.InsertAfter("Einleitung")
$ip3.LanguageID = "de-DE"
$ip3.Collapse(0)

# ---------------------------------------------------------------------------
# 4. Paragraph-format fix: the lone empty paragraph (right after the kernel-
#    size bullet list, right before the "Aufbau des neuronalen Netzes"
#    heading) that had explicit justification now instead carries the
#    Textkoerper style (with the green accent3 color preserved in its run
#    properties).
# ---------------------------------------------------------------------------
$fixIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count - 1; $i++) {
    $p = $d.Paragraphs.Item($i)
    $next = $d.Paragraphs.Item($i + 1)
    if ($p.Style.NameLocal -eq "Normal" -and $p.Alignment -eq 3 -and $next.Range.Text -like "Aufbau des neuronalen Netzes*") {
        $fixIdx = $i
        break
    }
}

if ($fixIdx -ne $null) {
    $fixPara = $d.Paragraphs.Item($fixIdx)
    $fixRange = $fixPara.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes" ?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:color w:val="9BBB59" w:themeColor="accent3"/><w:lang w:val="de-DE"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $fixRange.InsertXML($xml)
}

Write-Output "done"
